$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 0.8213140000000001
$ws.Range("N2").Value = 2.463942
$ws.Range("O2").Value = 0.06824749762056036
$ws.Range("P2").Value = 0.06824749762056037
$ws.Range("Q2").Value = 0.3748509948560001
$ws.Range("R2").Value = 3.373658953704
$ws.Range("S2").Value = 0.001306912803896207
$ws.Range("T2").Value = 0.001306912803896208
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.2017018900182306
$ws.Range("P3").Value = 0.2017018900182306
$ws.Range("Q3").Value = 1.107852401534667
$ws.Range("R3").Value = 9.970671613812002
$ws.Range("S3").Value = 0.003862512060156116
$ws.Range("T3").Value = 0.003862512060156116
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.7300506123612091
$ws.Range("P4").Value = 0.7300506123612091
$ws.Range("Q4").Value = 4.009820255393334
$ws.Range("R4").Value = 36.08838229854
$ws.Range("S4").Value = 0.01398018280599483
$ws.Range("T4").Value = 0.01398018280599483
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 0.8213140000000001
$ws.Range("N5").Value = 2.463942
$ws.Range("O5").Value = 0.06824749762056036
$ws.Range("P5").Value = 0.06824749762056037
$ws.Range("Q5").Value = 16.21782421005934
$ws.Range("R5").Value = 145.960417890534
$ws.Range("S5").Value = 0.05654321957877333
$ws.Range("T5").Value = 0.05654321957877334
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.2017018900182306
$ws.Range("P6").Value = 0.2017018900182306
$ws.Range("S6").Value = 0.1671105118045903
$ws.Range("T6").Value = 0.1671105118045903
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.7300506123612091
$ws.Range("P7").Value = 0.7300506123612091
$ws.Range("S7").Value = 0.6048487273168807
$ws.Range("T7").Value = 0.6048487273168807
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 0.8213140000000001
$ws.Range("N8").Value = 2.463942
$ws.Range("O8").Value = 0.06824749762056036
$ws.Range("P8").Value = 0.06824749762056037
$ws.Range("Q8").Value = 2.982190312686001
$ws.Range("R8").Value = 26.83971281417401
$ws.Range("S8").Value = 0.01039736523789083
$ws.Range("T8").Value = 0.01039736523789083
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.2017018900182306
$ws.Range("P9").Value = 0.2017018900182306
$ws.Range("S9").Value = 0.03072886615348419
$ws.Range("T9").Value = 0.03072886615348419
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.7300506123612091
$ws.Range("P10").Value = 0.7300506123612091
$ws.Range("S10").Value = 0.1112217022383336
$ws.Range("T10").Value = 0.1112217022383336
